# Daily attendance processing - 2026-01-16 11:07:16
# Swap the order of "dnasr281@gmail.com" and "System" in the
# "Recorded By" column (column G) wherever both appear together.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = 7 ("Recorded By")
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
    }
}
